# Apply cryptocurrency price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.235.06"
$ws.Range("E2").Value = "  +3.20%  "
$ws.Range("D3").Value = "2.064.94"
$ws.Range("E3").Value = "  +2.42%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'230.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("D7").Value = "'61.19"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +10.12%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +3.23%  "
$ws.Range("E10").Value = "  +3.91%  "
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("D12").Value = "'14.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.68%  "
$ws.Range("D13").Value = "2.370.60"
$ws.Range("E13").Value = "  +2.43%  "
$ws.Range("D14").Value = "'21.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.88%  "
$ws.Range("D15").Value = "'0.764"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.70%  "
$ws.Range("D16").Value = "'5.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.72%  "
$ws.Range("D17").Value = "2.082.05"
$ws.Range("E17").Value = "  +3.26%  "
$ws.Range("D18").Value = "38.096.02"
$ws.Range("E18").Value = "  +2.99%  "
$ws.Range("D19").Value = "'6.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.73%  "
$ws.Range("D20").Value = "'70.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("E21").Value = "  +3.12%  "
$ws.Range("D22").Value = "'226.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").Value = "'2.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.20%  "
$ws.Range("D26").Value = "'166.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("D27").Value = "'9.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("D28").Value = "'0.133"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.53%  "
$ws.Range("D29").Value = "'18.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.49%  "
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("D32").Value = "'4.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.06%  "
$ws.Range("D33").Value = "'4.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.52%  "
$ws.Range("E34").Value = "  +9.57%  "
$ws.Range("D35").Value = "'0.0605"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("B36").Value = "THORChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D36").Value = "'6.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +14.69%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("D38").Value = "'3.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.01%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "1.521.58"
$ws.Range("E40").Value = "  +3.56%  "
$ws.Range("D41").Value = "'98.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.12%  "
$ws.Range("D42").Value = "'17.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.02%  "
$ws.Range("E43").Value = "  +2.17%  "
$ws.Range("D44").Value = "'2.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.85%  "
$ws.Range("D45").Value = "'0.0926"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.11%  "
$ws.Range("D46").Value = "'1.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("E48").Value = "  +2.26%  "
$ws.Range("D49").Value = "'2.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("D50").Value = "'7.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("D51").Value = "2.257.94"
$ws.Range("E51").Value = "  +2.41%  "
